$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: H1 "syst_u" -> "syst_c"
$ws.Range("H1").Value = "syst_c"

# Data rows 2-9: boson column (E) "W" -> "W+"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = "W+"
}

# Re-saved numeric values carry the same magnitude but a different last-bit
# floating point representation (as Excel itself recomputed/re-serialized
# them) - reproduce exactly so the stored doubles match bit-for-bit.
$ws.Cells.Item(2, 6).Value = "0.26739999999999997"
$ws.Cells.Item(2, 8).Value = "0.008199999999999999"
$ws.Cells.Item(3, 7).Value = "0.0040999999999999995"
$ws.Cells.Item(5, 6).Value = "0.16260000000000002"
$ws.Cells.Item(6, 8).Value = "0.009000000000000001"
$ws.Cells.Item(7, 7).Value = "0.005699999999999999"
$ws.Cells.Item(7, 8).Value = "0.013500000000000002"
$ws.Cells.Item(9, 7).Value = "0.012199999999999999"
$ws.Cells.Item(9, 8).Value = "0.013000000000000001"

# Update the active selection to J21
$ws.Range("J21").Select()
